# "made type split a little more complex"
#
# The dwelling-type split model now phases the new type in gradually
# instead of having it appear at full share from the start:
#   - years 1600-1900 (rows 2-302): this type contributes 0 dwellings
#   - years 1901-2050 (rows 303-452): this type ramps up from a small
#     share toward the full total over time
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Years 1600-1900 -> column B (AB) is now 0 for this type.
for ($row = 2; $row -le 302; $row++) {
    $ws.Cells.Item($row, 2).Value = 0
}

# Years 1901-2050 -> ramped-up values for this type.
$rampValues = @(
    1952.789333408346,
    3936.550948613783,
    5951.819019872233,
    7999.134737616318,
    10079.04615512657,
    12192.10801445265,
    14338.88155074149,
    16519.93427413138,
    18735.83972837708,
    20987.17722538473,
    23274.53155484689,
    25598.49266818397,
    27959.65533602314,
    30358.61877846942,
    32795.98626745229,
    35272.36470046929,
    37788.3641450828,
    40344.59735357259,
    42941.67924719788,
    45580.22636957189,
    48260.85630871757,
    50984.18708743298,
    53750.83652167212,
    56561.42154671764,
    59416.55751100808,
    62316.85743757067,
    65262.9312531027,
    68255.38498484596,
    71294.81992550446,
    74381.83176656617,
    77517.00970050534,
    80700.93549246514,
    83934.18252214605,
    87217.31479675618,
    90550.88593601692,
    93935.43813035452,
    97371.50107355112,
    100859.5908712779,
    104400.2089270763,
    107993.8408075043,
    111640.9550883193,
    115342.0021837146,
    119097.4131607829,
    122907.5985415227,
    126772.9470948619,
    130693.8246213045,
    134670.5727329602,
    138703.5076318388,
    142792.9188894376,
    146939.0682307571,
    151142.1883260095,
    155402.4815933823,
    159720.1190163204,
    164095.238978878,
    168527.9461227604,
    173631.81597256,
    178801.5169093511,
    184037.0649203869,
    189338.4338358094,
    194705.5542574598,
    200138.3125385466,
    205636.5498183593,
    211200.0611161471,
    216828.5944881952,
    222521.8502520359,
    228279.4802815879,
    234101.0873768837,
    239986.2247118619,
    245934.3953635148,
    251945.0519254718,
    257068.0633228233,
    262242.0257195677,
    267466.3003985764,
    272740.2032940052,
    278063.0049828968,
    283433.930758112,
    288852.1607835802,
    294316.8303325808,
    299827.030109481,
    305381.8066550634,
    309427.5267279205,
    313501.8310327973,
    317603.8043896527,
    321732.4960398872,
    325886.9200870801,
    330066.0559928642,
    334268.8491267727,
    338494.2113687384,
    342741.0217627939,
    347008.127220402,
    354134.9439817925,
    361299.9953418839,
    368501.6846557999,
    375738.3842124285,
    383008.437146236,
    390310.1594123216,
    397641.8418204443,
    405001.7521236516,
    412388.1371570933,
    419799.2250225376,
    431434.2978306849,
    443113.8589984321,
    454835.6417335415,
    466597.381637312,
    478396.8206347519,
    490231.7109218517,
    502099.8189209187,
    513998.9292347442,
    525926.8485902692,
    537881.4097621809,
    549257.0474527335,
    560652.4093320038,
    572065.4650932846,
    583494.2220171342,
    594936.7284535042,
    606391.077227426,
    617855.4089570551,
    629327.9152723069,
    640806.8419216803,
    652290.4917541472,
    663777.227562227,
    675265.4747714817,
    686753.7239608071,
    698240.5331968718,
    709724.5301650916,
    721204.414078468,
    732678.9573445545,
    744147.0069697828,
    755607.4856793466,
    767059.3927298723,
    778501.8043912349,
    789933.8740730986,
    801354.8320711723,
    812763.9849077208,
    824160.7142407335,
    835544.4753161991,
    846914.7949383962,
    858271.2689338492,
    869613.5590858196,
    880941.3895178344,
    892254.5425068727,
    903552.8537095218,
    914836.206787614,
    926104.5274236699,
    937357.7767208819,
    948595.9439873892,
    959819.0389102506,
    971027.0831307045,
    982220.1012391638,
    993398.1112157055
)

$firstRampRow = 303
for ($i = 0; $i -lt $rampValues.Length; $i++) {
    $ws.Cells.Item($firstRampRow + $i, 2).Value = $rampValues[$i]
}

Write-Host "zeroed $(302 - 2 + 1) cells and ramped $($rampValues.Length) cells in column B"
